# Generate Report for Handoff
# Updates the "b.md" rows across the Overview, zh-cn and de-de sheets to
# reflect that a new handoff package was generated for b.md.

$wb = $excel.ActiveWorkbook

$newStatus      = "Ready for handoff"
$newStatusDate  = "2017-02-17 07:13:21"

# ---------------------------------------------------------------------
# Overview sheet - row for b.md (row 3)
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus
$wsOverview.Range("G3").Value = $newStatusDate

# ---------------------------------------------------------------------
# zh-cn sheet - row for b.md (row 3)
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = $newStatus
# Leading apostrophe forces literal text (avoids "False" being auto-coerced
# to a Boolean); reset Style afterwards so no stray quote-prefix formatting
# is left behind on the cell.
$wsZhCn.Range("F3").Value = "'False"
$wsZhCn.Range("F3").Style = "Normal"
$wsZhCn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2017-02-17 07:13:04"
$wsZhCn.Range("R3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/769d532d8485574c5025c086d82bfc29326a3c70/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2775b4392889070fb7da6e2a171ebc90a3eb3114/e2e/b.md."
$wsZhCn.Columns.Item(18).ColumnWidth = 39.2

# ---------------------------------------------------------------------
# de-de sheet - row for b.md (row 3)
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = $newStatus
$wsDeDe.Range("F3").Value = "'False"
$wsDeDe.Range("F3").Style = "Normal"
$wsDeDe.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDeDe.Range("H3").Value = $newStatusDate
$wsDeDe.Range("R3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/769d532d8485574c5025c086d82bfc29326a3c70/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2775b4392889070fb7da6e2a171ebc90a3eb3114/e2e/b.md."
$wsDeDe.Columns.Item(18).ColumnWidth = 39.2
